$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Number" column values
$ws.Range("E2").Value = 3
$ws.Range("E3").Value = 6
$ws.Range("E4").Value = 6
$ws.Range("E5").Value = 10

# Update the selected cell/range to E5
$ws.Range("E5").Select()
